$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before D, shifting old D:K to F:M
$ws.Range("D:E").EntireColumn.Insert()

# Copy number formats/styles from F:G (shifted originals) into the new D:E columns
$ws.Range("F:G").Copy()
$ws.Range("D:E").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Populate new column D and E values for the new quarters (2018-12-31 and 2018-09-30)
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 589600
$ws.Range("E8").Value = 302500
$ws.Range("D9").Value = 480700
$ws.Range("E9").Value = 316100
$ws.Range("D10").Value = 108900
$ws.Range("E10").Value = -13600
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = -9300
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 24900
$ws.Range("E15").Value = 22400
$ws.Range("D17").Value = 499600
$ws.Range("E17").Value = 341100
$ws.Range("D18").Value = 90000
$ws.Range("E18").Value = -38600
$ws.Range("D20").Value = -2800
$ws.Range("E20").Value = 1400
$ws.Range("D21").Value = 89900
$ws.Range("E21").Value = -39600
$ws.Range("D22").Value = 30200
$ws.Range("E22").Value = 26500
$ws.Range("D23").Value = 57000
$ws.Range("E23").Value = -63700
$ws.Range("D24").Value = 12800
$ws.Range("E24").Value = -16600
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 44300
$ws.Range("E26").Value = -47100
$ws.Range("D27").Value = 46000
$ws.Range("E27").Value = -45600
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = -100
$ws.Range("E29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 2800
$ws.Range("E32").Value = -1400
$ws.Range("D33").Value = 45900
$ws.Range("E33").Value = -45700
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 45900
$ws.Range("E35").Value = -45700
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 30000
$ws.Range("E41").Value = 3300
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 320600
$ws.Range("E43").Value = 203800
$ws.Range("D44").Value = 62200
$ws.Range("E44").Value = 88900
$ws.Range("D45").Value = 250400
$ws.Range("E45").Value = 441200
$ws.Range("D46").Value = 663200
$ws.Range("E46").Value = 737300
$ws.Range("D47").Value = 117400
$ws.Range("E47").Value = 116100
$ws.Range("D48").Value = 3653500
$ws.Range("E48").Value = 3530100
$ws.Range("D49").Value = 734600
$ws.Range("E49").Value = 759800
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 787900
$ws.Range("E52").Value = 771200
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 5956600
$ws.Range("E54").Value = 5914500
$ws.Range("D57").Value = 410500
$ws.Range("E57").Value = 383500
$ws.Range("D58").Value = 1004400
$ws.Range("E58").Value = 1904100
$ws.Range("D59").Value = 166000
$ws.Range("E59").Value = 183600
$ws.Range("D60").Value = 1580800
$ws.Range("E60").Value = 2471200
$ws.Range("D61").Value = 2106900
$ws.Range("E61").Value = 1281000
$ws.Range("D62").Value = 1001900
$ws.Range("E62").Value = 927500
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 4689600
$ws.Range("E66").Value = 4679700
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 343300
$ws.Range("E72").Value = 321900
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 1267000
$ws.Range("E76").Value = 1234800
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 45900
$ws.Range("E81").Value = -45700
$ws.Range("D83").Value = 0
$ws.Range("E83").Value = 0
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = -29800
$ws.Range("E89").Value = 20600
$ws.Range("D91").Value = -140400
$ws.Range("E91").Value = -74800
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = 171800
$ws.Range("E94").Value = -1828700
$ws.Range("D96").Value = -48500
$ws.Range("E96").Value = -23900
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -123600
$ws.Range("E100").Value = 48000
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 18400
$ws.Range("E102").Value = -1760100
